$titles = @(
    'Alienware Headset Inalámbrico, Blanco',
    'Bulova Winston C4844 - Reloj de Pared (1 Unidad), Color Plateado.',
    'TOPLIVING Silla de Oficina Tipo Gerencial con Cabecera Ajustable Base Metálica para Home Office Silla, con Base Giratoria, Altura Regulable, Color Gris',
    'Simple Wisdom Librero de 5 Capas, Librero Organizador, Trapezoidal Librero Minimalista, Estante Organizador de Uso Multiescenario Adecuado para Sala de Estar, Oficina, Dormitorio, Cocina, Baño（Negro）',
    'Tent SUNDOME',
    'Lenovo Laptop Gaming LOQ 15IAX9E 15.6 Pulgadas Intel Core i5-12450HX 8GB 512 GB SSD RTX3050 de 8 núcleos, Windows 11 Home Color Luna Grey',
    'Bulova Clocks Precise Connect - Reloj de Pared C5001, 12.5 Pulgadas, Plata Mate y Azul',
    'UGREEN Cable USB Micro Nylon Trenzado Rápido Cargador Compatible con Galaxy S7 S6 Edge S5 Nota 5, Xiaomi, Huawei, HTC, Oneplus, LG, E lectores, 1M',
    'Hyperkin Xenon Wired Controller (Black) For Xbox Series X|S/Xbox One/Windows 10|11',
    'Amazon Basics - Juego de sábanas de microfibra ligera y súper suave de fácil cuidado con bolsillos de 14 pulgadas de profundidad, individual, color crema',
    'CRAFTSMAN Batería de Iones de Litio V20, 4.0 amperios Hora (CMCB204)',
    'Bulova Clocks C4834 Luxfer Prism - Reloj de Pared Inspirado en Frank Lloyd Wright de 14 Pulgadas',
    'Ernie Ball Earthwood Medium - Cuerdas para guitarra acústica de 12 cuerdas, bronce 80/20, calibre 11-28',
    'Anker 20W Hielo Cargador Rápido, Control de Temperatura Avanzado, Diseño Compacto 25% Más Pequeño, Compatible con iPhone y Samsung y Otros Dispositivos, USB C Tipo, 2 Colores',
    'EASYTAO Tira Led Neón 5M 600 LED, IP65 Impermeable, Luces LED Flexible de Decoración para Adecuado TV, Salón, Dormitorio, Cuarto (Blanco Frío)',
    'EASYTAO TV-308 Cámara de Seguridad Interior 2.5K, 4MP Cámara Alexa WiFi, iCSee y Tris Home App, Audio Bidireccional, Control Remoto, Detección de Movimiento, Visión Nocturna,Compatible con Alexa',
    'Logitech G Driving Force Shifter Palanca de Cambios para Volantes de Carreras G29, G920 y G923, 6 Velocidades, Marcha Atrás a Presion, Acero y Cuero - Negro',
    'Fabuloso Limpiador Multiusos Líquido Frescura Activa Lavanda 4 L. Antiviral y Antibacterial. Limpiador Desinfectante que Brinda 24 Horas de Frescura',
    'Fancy Pets Mueble/Rascador para Gato Tebas de 140 Centímetros de Altura',
    'BYASKZIV Cargador Tipo C (1.5M)',
    'Bluelander Mochila Cruzada Antirrobo Impermeable, Bolso de Pecho para Hombre, Bolso Bandolera con Puerto de Carga USB, Perfecta para Celular, Tablets, Actividades como Deportes, Compras',
    'Raganet, Mochila Transportadora de Mascotas, Perros y Gatos, Back Pack Portátil para Mascota Pequeña (Máximo 7Kg) Medidas: Ancho 23cm Largo 33cm Altura 42cm (Color Negro)',
    'Nescafé Dolce Gusto Café en Cápsulas Diferentes Variedades 24 Cápsulas',
    'Feandrea Árbol para Gatos de 155 cm de Altura, 49 x 40 x 155 cm, Rascador para Gatos con 5 Postes, Felpa Multicapa, Blanco Crema MPCT192T01',
    'HIKEO Hielera Térmica Extra Grande Impermeable (30 Litros/36 Latas), Mochila Térmica Portátil con Múltiples Compartimentos, Doble Pared a Prueba de Fugas para Hielo, Ideal para Camping, Picnic, Viajes',
    'KitKat chocolate con leche 9 pz de 41.5gr',
    'JBL Bar 300 Barra de Sonido 5.0 Wi-Fi/Bluetooth, 260W de Potencia, Dolby Atmos, Multibeam, con HDMI ARC - Negro',
    'Cama Matrimonial Base Y Cabecero de Almacenamiento con Luces Led 2 Salidas de CA y 2 USB,Con dos cajones extraíbles y Cabecero Plataforma Tapizada,Antideslizante,sin Ruido,Café Negro(Matrimonial/Full)',
    'Tira Led Neon RGB 30M Luces Led Para Cuarto Exterior Flex Manguera Party Luz Lights Navidad Halloween Fiesta Lampara Habitacion Aesthetic Calida Leds Strip Home Decoracion 2 Rollos',
    'Deyidn Lámpara de tocador moderna de 2 Luces, oro lámparas Vidrio acrílico espejo LED luz ajustable moderna lámpara de pared escocesa (12W 4000k)',
    'Bluelander Radios de Comunicación Inalámbricos 2 Piezas, Walkie Talkie de Largo Alcance, Walkie Talkie Recargable con Base de Carga a Prueba de Agua, Radios Walkie Talkie para Guardia, Almacé, Camping',
    'CYCROWN CycRun Bicicleta Electrica para Adultos, Moto Electrica con Potencia Máxima del Motor 1000W,Autonomía 70-90KM, Velocidad de 35KM/H, con 20" Neumático Moto Bicicleta de Montaña',
    'MURAMASA Mesa De Centro para Sala Mesa De Centro Elevable con Compartimentos Ocultos y Estante De Almacenamiento Mueble para Sala con Estación De Carga Mesa Expandible Madera Negro',
    'Jimmy Choo Blossom, 3.3 Fluid Ounce',
    'Twinings, Twinings Pure Green Tea, 25 sobres de 2 gramos cada uno',
    '120cm Credenza Recibidor con Tomacorrientes y Puertos USB,Mesa Consola de Entrada con 3 cajones y 3 Estantes,Muebles Mesa Consola para Sala,Recámara,Comedor,Soporte de TV,Barra de Café,Mesa de Sofá',
    'Funko Pop! Deluxe: Cuphead - The Devil In Chair - Figura de Vinilo Coleccionable - Idea de Regalo - Mercancia Oficial - Juguetes para Niños y Adultos - Video Games Fans - Muñeco para Coleccionistas',
    'amiibo - Cat Mario - Super Mario Series',
    'NESCAFÉ Dolce Gusto Cappuccino Skinny Light de 16 cápsulas',
    'TURTLE BEACH RECON 70 BLUE CAMO MULTIPLATAFORMA'
)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3
for ($i = 0; $i -lt $titles.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $titles[$i]
}
